$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the prompt text in C5 (Effectively deploy LLMs task)
$ws.Range("C5").Value = "How to most effectively deploy LLMs?"

# Update the active selection to reflect the last edited cell
[void]$ws.Range("C5").Select()
